# This workbook holds a rolling window of trading-day data in three
# worksheets (columns A = date serial, B = remn_amt), each with rows 2..101
# holding 100 data rows below the header in row 1.
#
# The update "rolls" the window forward by two trading days:
#   - the two oldest rows (old rows 2 and 3) are dropped
#   - every remaining row shifts up by two positions (old row N -> new row N-2)
#   - two brand new rows are appended at the bottom with the next two
#     calendar dates and a placeholder remn_amt of 0 (data not yet available)

$wb = $excel.ActiveWorkbook

$firstDataRow = 2
$lastDataRow = 101
$numRows = $lastDataRow - $firstDataRow + 1   # 100
$shift = 2

for ($s = 1; $s -le $wb.Worksheets.Count; $s++) {
    $ws = $wb.Worksheets.Item($s)

    $readRange = $ws.Range("A$firstDataRow" + ":B$lastDataRow")
    $orig = $readRange.Value2

    $new = New-Object 'object[,]' $numRows,2

    # Shift every row up by $shift positions: new row i gets the values
    # that used to live $shift rows further down.
    for ($i = 1; $i -le ($numRows - $shift); $i++) {
        $new[$i - 1, 0] = $orig[$i + $shift, 1]
        $new[$i - 1, 1] = $orig[$i + $shift, 2]
    }

    # Append $shift new rows at the end: dates continue the existing
    # daily sequence, remn_amt starts at 0 until real data arrives.
    $lastDate = $orig[$numRows, 1]
    for ($k = 1; $k -le $shift; $k++) {
        $rowIdx = $numRows - $shift + $k
        $new[$rowIdx - 1, 0] = $lastDate + $k
        $new[$rowIdx - 1, 1] = 0
    }

    $readRange.Value2 = $new
}
